$p = $ppt.ActivePresentation
$s1 = $p.Slides.Item(1)
Write-Host "Before:" $s1.ColorScheme.Item(3).RGB
try {
  $s1.ColorScheme = $s1.ColorScheme
  Write-Host "assign self ok"
} catch {
  Write-Host "ERR:" $_.Exception.Message
}
